# dsa arrays and math
# Adds two new rows (86, 87) to Sheet1 describing two new LeetCode problems:
#   28. Find the Index of the First Occurrence in a String (Arrays)
#   69. Sqrt(x) (Math)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 86: "28. Find the Index of the First Occurrence in a String" ---
$ws.Range("A86").Value = '28. Find the Index of the First Occurrence in a String'
$ws.Range("B86").Value = 'Easy'
$ws.Range("B86").Interior.Color = 5287936
$ws.Range("C86").Value = 'Arrays'
$ws.Range("D86").Value = 'We can use a 2d for loop to iterate the haystack''s chars and needle''s chars. If i + needle length > haystack length, immediately break, iterate if haystack.charAt(i+j) is equal to needle.charAt(j) and iterate until j == needle length. '
$ws.Range("E86").Value = 'https://leetcode.com/problems/find-the-index-of-the-first-occurrence-in-a-string/solutions/12807/elegant-java-solution/comments/142986 '
$ws.Hyperlinks.Add($ws.Range("E86"), 'https://leetcode.com/problems/find-the-index-of-the-first-occurrence-in-a-string/solutions/12807/elegant-java-solution/comments/142986 ')
$ws.Range("E86").Style = "Hyperlink"

# --- Row 87: "69. Sqrt(x)" ---
$ws.Range("A87").Value = '69. Sqrt(x)'
$ws.Range("B87").Value = 'Easy'
$ws.Range("B87").Interior.Color = 5287936
$ws.Range("C87").Value = 'Math'
$ws.Range("E87").Value = 'https://leetcode.com/problems/sqrtx/solutions/25057/3-4-short-lines-integer-newton-every-language/?envType=study-plan-v2&envId=top-interview-150 '
$ws.Hyperlinks.Add($ws.Range("E87"), 'https://leetcode.com/problems/sqrtx/solutions/25057/3-4-short-lines-integer-newton-every-language/?envType=study-plan-v2&envId=top-interview-150 ')
$ws.Range("E87").Style = "Hyperlink"
$ws.Range("D87").Value = 'Optimally, use the Newton method. Otherwise, use binary search with the search condition of the square of mid > mid.'

# Restore the selection to roughly where the author left off.
$ws.Range("D92").Select()
